# Edit script: adds "PO Forecast" sheet, renames headers on existing sheets,
# and populates forecast data.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Requested quantity" headers on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" worksheet as the last sheet in the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- 3. Header row ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# --- 4. Data rows (ds, PO_Forecast, yhat_lower, yhat_upper) ---
$data = @(
    @(44948.99999999999,47,-18.07049198609675,111.3428414981349),
    @(44976.99999999999,47,-19.24527381582851,119.0144614903215),
    @(44983.99999999999,47,-20.38128186169032,111.3598371895123),
    @(44990.99999999999,47,-19.07379848721462,118.1604287254523),
    @(45011.99999999999,47,-22.40831867756932,112.104254657724),
    @(45032.99999999999,47,-20.62923113842756,112.8207925625254),
    @(45039.99999999999,47,-19.36384117982274,117.4541956291761),
    @(45102.99999999999,47,-15.22431122566123,112.9749838837237),
    @(45179.99999999999,47,-16.28020446137431,112.9955122723683),
    @(45186.99999999999,47,-21.36648739283071,112.9691443264765),
    @(45193.99999999999,48,-21.10087178760246,113.2934721942145),
    @(45200.99999999999,48,-19.53877714539983,117.3035849944991),
    @(45207.99999999999,48,-21.54459920934859,116.0314822259014),
    @(45214.99999999999,48,-18.91583763104572,116.4016640807074),
    @(45221.99999999999,48,-19.51217848392528,115.1511471406702),
    @(45228.99999999999,48,-20.66723379346599,113.2928776663304),
    @(45263.99999999999,48,-22.90982839859288,117.1613181215433),
    @(45340.99999999999,48,-22.14729841197099,114.3063450413549),
    @(45347.99999999999,48,-17.96532671255548,117.393934248601),
    @(45354.99999999999,48,-24.30628948717745,113.5593750687312),
    @(45361.99999999999,48,-18.14632428375663,114.1529445209453),
    @(45382.99999999999,48,-20.68386566809389,110.7590449928255),
    @(45389.99999999999,48,-19.25386550696296,115.1145973734048),
    @(45431.99999999999,48,-16.83009153234638,116.9323621455069),
    @(45466.99999999999,48,-19.95989746361629,117.0286302746869),
    @(45480.99999999999,48,-14.37427114508775,114.9458378290959),
    @(45536.99999999999,48,-16.79749708897091,115.2768042445973),
    @(45592.99999999999,48,-19.44239367530855,116.2163623988825),
    @(45599.99999999999,48,-13.25978375999855,116.7428536400458),
    @(45606.99999999999,48,-18.32050287123608,118.9047025069304),
    @(45613.99999999999,48,-21.77540561768442,116.5191499422867),
    @(45620.99999999999,48,-13.21629307067905,114.840656727519),
    @(45627.99999999999,48,-21.74286045925234,112.8574079146694),
    @(45634.99999999999,48,-23.97832202566746,117.7060150468031),
    @(45641.99999999999,48,-18.36657786532348,116.8514272127114),
    @(45648.99999999999,48,-19.27251698956994,116.9227269324218),
    @(45655.99999999999,48,-17.65447990993656,118.5386706739358)
)

$rowCount = $data.Count
$colCount = 4
$values = New-Object 'object[,]' $rowCount, $colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    for ($j = 0; $j -lt $colCount; $j++) {
        $values[$i, $j] = $data[$i][$j]
    }
}

$startRow = 2
$endRow = $startRow + $rowCount - 1
$dataRangeAddress = "A${startRow}:D${endRow}"
$dataRange = $wsForecast.Range($dataRangeAddress)
$dataRange.Value = $values

# Apply the same date/time number format used for the "ds" columns on the other sheets
$dateColAddress = "A${startRow}:A${endRow}"
$wsForecast.Range($dateColAddress).NumberFormat = "YYYY-MM-DD HH:MM:SS"
